$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. Remove the floating "MEDIA: Media de idade na tabela" text-box shape
# ------------------------------------------------------------------
if ($ws.Shapes.Count -gt 0) {
    for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
        $ws.Shapes.Item($i).Delete()
    }
}

# ------------------------------------------------------------------
# 2. Unmerge the D:E label/value blocks (MEDIA DE IDADE, MEDIANA, VALOR DO MEIO)
# ------------------------------------------------------------------
$ws.Range("D6:E6").UnMerge()
$ws.Range("D9:E9").UnMerge()
$ws.Range("D10:E10").UnMerge()

# ------------------------------------------------------------------
# 3. Drop the "ANOS" unit labels next to the average/median results,
#    leaving the cells blank (keeping their existing formatting)
# ------------------------------------------------------------------
$ws.Range("E7").ClearContents()
$ws.Range("E11").ClearContents()

# ------------------------------------------------------------------
# 4. Re-style the now-unmerged D/E cells: a full box border around the
#    value cells in column D, no border / no horizontal centering on
#    the now-empty column E cells.
# ------------------------------------------------------------------
$ws.Range("D6").Borders.LineStyle = 1
$ws.Range("D6").Borders.Weight = 2
$ws.Range("D6").HorizontalAlignment = -4108
$ws.Range("D6").VerticalAlignment = -4108

$ws.Range("E6").Borders.LineStyle = 0
$ws.Range("E6").VerticalAlignment = -4108

$ws.Range("E7").Borders.LineStyle = 0

$ws.Range("D9").Borders.LineStyle = 1
$ws.Range("D9").Borders.Weight = 2
$ws.Range("E9").Borders.LineStyle = 0
$ws.Range("E9").HorizontalAlignment = -4142

$ws.Range("D10").Borders.LineStyle = 1
$ws.Range("D10").Borders.Weight = 2
$ws.Range("D10").Font.Italic = $true
$ws.Range("E10").Borders.LineStyle = 0
$ws.Range("E10").HorizontalAlignment = -4142

$ws.Range("D11").Font.Bold = $true
$ws.Range("E11").Borders.LineStyle = 0

# New blank (bordered) cells under the average value
$ws.Range("D8").Value = ""
$ws.Range("E8").Value = ""

# Stray underline formatting left on a hidden helper cell
$ws.Range("D16").Font.Underline = $true

# ------------------------------------------------------------------
# 5. Hide the leftover helper rows below the data (15-26) and drop the
#    now-unused blank spacer row 14
# ------------------------------------------------------------------
$ws.Rows("14:14").Delete()
$ws.Rows("15:26").Hidden = $true

# ------------------------------------------------------------------
# 6. Column D now holds its own label text (no longer spanning into E);
#    size it to fit its longest entry ("MEDIA DE IDADE").
# ------------------------------------------------------------------
$ws.Columns("D:D").ColumnWidth = 14.7109375

# ------------------------------------------------------------------
# 7. Selection left on column F (mirrors the author clicking the column
#    header while tidying up the hidden helper columns)
# ------------------------------------------------------------------
$ws.Range("F1:G1048576").Select()
